$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")

# Row 48
$ws.Range("H48").Value = 5000
$ws.Range("J48").Value = 5000
$ws.Range("L48").Value = 15000
$ws.Range("N48").Value = -15584

# Row 56
$ws.Range("H56").Value = 5000
$ws.Range("J56").Value = 5000
$ws.Range("L56").Value = 15000
$ws.Range("N56").Value = -16068

# Row 87
$ws.Range("H87").Value = 24500
$ws.Range("J87").Value = 24500
$ws.Range("L87").Value = 24500
$ws.Range("N87").Value = -26996

# Row 90
$ws.Range("H90").Value = 24500
$ws.Range("J90").Value = 24500
$ws.Range("L90").Value = 73500
$ws.Range("N90").Value = -85980

# Row 96
$ws.Range("H96").Value = 425.93332
$ws.Range("I96").Value = 424.15384
$ws.Range("K96").Value = 1272.46152
$ws.Range("M96").Value = 100.5384799999999

# Row 112
$ws.Range("H112").Value = 2007
$ws.Range("I112").Value = 1094.25
$ws.Range("K112").Value = 3282.75
$ws.Range("M112").Value = -2174.75

# Row 113
$ws.Range("H113").Value = 6943.778
$ws.Range("I113").Value = 5833.3335
$ws.Range("K113").Value = 5833.3335
$ws.Range("M113").Value = -2579.3335

# Row 116
$ws.Range("H116").Value = 4442.6665
$ws.Range("J116").Value = 4442.6665
$ws.Range("L116").Value = 4442.6665
$ws.Range("N116").Value = -11326.6665

# Row 137
$ws.Range("H137").Value = 3541.2222
$ws.Range("I137").Value = 2553.1428
$ws.Range("J137").Value = 6999.5
$ws.Range("K137").Value = 7659.428400000001
$ws.Range("L137").Value = 20998.5
$ws.Range("M137").Value = -5109.428400000001
$ws.Range("N137").Value = -26098.5

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")

# Row 3
$ws.Range("H3").Value = 1250
$ws.Range("I3").Value = 500
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 500
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -385
$ws.Range("N3").Value = -2230

# Row 5
$ws.Range("H5").Value = 568.375
$ws.Range("I5").Value = 752.5
$ws.Range("J5").Value = 384.25
$ws.Range("K5").Value = 752.5
$ws.Range("L5").Value = 384.25
$ws.Range("M5").Value = -640.5
$ws.Range("N5").Value = -608.25

# Row 17
$ws.Range("H17").Value = 6750
$ws.Range("J17").Value = 6750
$ws.Range("L17").Value = 6750
$ws.Range("N17").Value = -7096

# Row 32
$ws.Range("H32").Value = 2847.3447
$ws.Range("I32").Value = 1203.4348
$ws.Range("J32").Value = 9149
$ws.Range("K32").Value = 1203.4348
$ws.Range("L32").Value = 9149
$ws.Range("M32").Value = -916.4348
$ws.Range("N32").Value = -9723

# Row 58
$ws.Range("H58").Value = 19350
$ws.Range("J58").Value = 19350
$ws.Range("L58").Value = 19350
$ws.Range("N58").Value = -20210

# Row 61
$ws.Range("H61").Value = 2349.4
$ws.Range("I61").Value = 1542.4286
$ws.Range("J61").Value = 4232.3335
$ws.Range("K61").Value = 1542.4286
$ws.Range("L61").Value = 4232.3335
$ws.Range("M61").Value = -1330.4286
$ws.Range("N61").Value = -4656.3335

# Row 80
$ws.Range("H80").Value = 30000
$ws.Range("J80").Value = 30000
$ws.Range("L80").Value = 30000
$ws.Range("N80").Value = -31996

# Row 83
$ws.Range("H83").Value = 30000
$ws.Range("J83").Value = 30000
$ws.Range("L83").Value = 90000
$ws.Range("N83").Value = -99984

# Row 102
$ws.Range("H102").Value = 1430
$ws.Range("I102").Value = 1430
$ws.Range("K102").Value = 1430
$ws.Range("M102").Value = 192

# Row 122
$ws.Range("H122").Value = 2012
$ws.Range("I122").Value = 2012
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6036
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3586
$ws.Range("N122").ClearContents()

# Row 136
$ws.Range("H136").Value = 2349.4
$ws.Range("I136").Value = 1542.4286
$ws.Range("J136").Value = 4232.3335
$ws.Range("K136").Value = 4627.2858
$ws.Range("L136").Value = 12697.0005
$ws.Range("M136").Value = -2077.2858
$ws.Range("N136").Value = -17797.0005

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")

# Row 4
$ws.Range("H4").Value = 568.375
$ws.Range("I4").Value = 752.5
$ws.Range("J4").Value = 384.25
$ws.Range("K4").Value = 752.5
$ws.Range("L4").Value = 384.25
$ws.Range("M4").Value = -637.5
$ws.Range("N4").Value = -614.25

# Row 7
$ws.Range("H7").Value = 616.7143
$ws.Range("I7").Value = 40
$ws.Range("J7").Value = 712.8333
$ws.Range("K7").Value = 40
$ws.Range("L7").Value = 712.8333
$ws.Range("M7").Value = 73
$ws.Range("N7").Value = -938.8333

# Row 10
$ws.Range("H10").Value = 1000
$ws.Range("J10").Value = 1000
$ws.Range("L10").Value = 1000
$ws.Range("N10").Value = -1280

# Row 12
$ws.Range("H12").Value = 4422
$ws.Range("I12").Value = 5152.5
$ws.Range("J12").Value = 1500
$ws.Range("K12").Value = 5152.5
$ws.Range("L12").Value = 1500
$ws.Range("M12").Value = -4984.5
$ws.Range("N12").Value = -1836

# Row 22
$ws.Range("H22").Value = 332.5
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

# Row 94
$ws.Range("H94").Value = 2246.6
$ws.Range("J94").Value = 1999
$ws.Range("L94").Value = 1999
$ws.Range("N94").Value = -2901

# Row 107
$ws.Range("H107").Value = 4379.5
$ws.Range("J107").Value = 4748.5
$ws.Range("L107").Value = 4748.5
$ws.Range("N107").Value = -8588.5

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")

# Row 5
$ws.Range("H5").Value = 233
$ws.Range("I5").Value = 179.6
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 179.6
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = -67.59999999999999
$ws.Range("N5").Value = -724

# Row 31
$ws.Range("H31").Value = 1251.1818
$ws.Range("J31").Value = 1158.6
$ws.Range("L31").Value = 1158.6
$ws.Range("N31").Value = -1748.6

# Row 34
$ws.Range("H34").Value = 1251.1818
$ws.Range("J34").Value = 1158.6
$ws.Range("L34").Value = 1158.6
$ws.Range("N34").Value = -1562.6

# Row 58
$ws.Range("H58").Value = 1571.8572
$ws.Range("I58").Value = 1565.6428
$ws.Range("J58").Value = 1584.2858
$ws.Range("K58").Value = 1565.6428
$ws.Range("L58").Value = 1584.2858
$ws.Range("M58").Value = -1362.6428
$ws.Range("N58").Value = -1990.2858

# Row 59
$ws.Range("H59").Value = 28626

# Row 107
$ws.Range("H107").Value = 347.73334
$ws.Range("I107").Value = 222.8
$ws.Range("J107").Value = 597.6
$ws.Range("K107").Value = 222.8
$ws.Range("L107").Value = 597.6
$ws.Range("M107").Value = 1697.2
$ws.Range("N107").Value = -4437.6

# Row 134
$ws.Range("H134").Value = 3074.4
$ws.Range("I134").Value = 2918
$ws.Range("J134").Value = 4482
$ws.Range("K134").Value = 8754
$ws.Range("L134").Value = 13446
$ws.Range("M134").Value = -6219
$ws.Range("N134").Value = -18516

# Row 136
$ws.Range("H136").Value = 1571.8572
$ws.Range("I136").Value = 1565.6428
$ws.Range("J136").Value = 1584.2858
$ws.Range("K136").Value = 4696.928400000001
$ws.Range("L136").Value = 4752.857400000001
$ws.Range("M136").Value = -2146.928400000001
$ws.Range("N136").Value = -9852.857400000001

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")

# Row 63
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

# Row 66
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

# Row 113
$ws.Range("H113").Value = 1066.8572
$ws.Range("J113").Value = 1286
$ws.Range("L113").Value = 3858
$ws.Range("N113").Value = -8198

# Row 122
$ws.Range("H122").Value = 349.6
$ws.Range("I122").Value = 312
$ws.Range("K122").Value = 2808
$ws.Range("M122").Value = -358

# Row 137
$ws.Range("H137").Value = 5000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 5000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 15000
$ws.Range("N137").Value = -25200
$ws.Range("M137").ClearContents()

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")

# Row 4
$ws.Range("H4").Value = 625
$ws.Range("J4").Value = 625
$ws.Range("L4").Value = 625
$ws.Range("N4").Value = -849

# Row 5
$ws.Range("H5").Value = 84800
$ws.Range("I5").Value = 84800
$ws.Range("K5").Value = 84800
$ws.Range("M5").Value = -84688

# Row 6
$ws.Range("H6").Value = 3752
$ws.Range("J6").Value = 4836
$ws.Range("L6").Value = 4836
$ws.Range("N6").Value = -5062

# Row 9
$ws.Range("H9").Value = 5335
$ws.Range("I9").Value = 570
$ws.Range("K9").Value = 570
$ws.Range("M9").Value = -400

# Row 16
$ws.Range("H16").Value = 3752
$ws.Range("J16").Value = 4836
$ws.Range("L16").Value = 4836
$ws.Range("N16").Value = -5336

# Row 27
$ws.Range("H27").Value = 6166.3335
$ws.Range("I27").Value = 3999
$ws.Range("J27").Value = 7250
$ws.Range("K27").Value = 3999
$ws.Range("L27").Value = 7250
$ws.Range("M27").Value = -3833
$ws.Range("N27").Value = -7582

# Row 104
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")

# Row 16
$ws.Range("H16").Value = 4500
$ws.Range("I16").Value = 4000
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 4000
$ws.Range("L16").Value = 5000
$ws.Range("M16").Value = -3830
$ws.Range("N16").Value = -5340

# Row 36
$ws.Range("H36").Value = 365000
$ws.Range("J36").Value = 365000
$ws.Range("L36").Value = 365000
$ws.Range("N36").Value = -366124

# Row 136
$ws.Range("H136").Value = 3359.8
$ws.Range("I136").Value = 2976.3845
$ws.Range("K136").Value = 8929.1535
$ws.Range("M136").Value = -6379.1535

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")

# Row 119
$ws.Range("H119").Value = 48499.75
$ws.Range("J119").Value = 48499.75
$ws.Range("L119").Value = 48499.75
$ws.Range("N119").Value = -58175.75
